# feat: add 2022-Q3 data
#
# The workbook has a "总计" (summary) sheet followed by one sheet per
# quarter (newest first). This adds a brand-new "2022-Q3" quarter:
#   1. a new worksheet "2022-Q3" is inserted right after "总计" (so it
#      becomes the second tab, pushing 2022-Q1/2021-Q4/... one slot right)
#   2. the "总计" sheet gets a new row 2 summarizing the new quarter, and
#      the existing summary rows shift down (their running index in
#      column A is renumbered to stay contiguous)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Insert the new "2022-Q3" worksheet right after "总计"
# ---------------------------------------------------------------------
$zongji = $wb.Worksheets.Item("总计")
$q3 = $wb.Worksheets.Add($null, $zongji)
$q3.Name = "2022-Q3"

# Header row (identical layout/text to the other quarter sheets)
$q3.Cells.Item(1, 2).Value = "基金代码"
$q3.Cells.Item(1, 3).Value = "基金名称"
$q3.Cells.Item(1, 4).Value = "基金规模"
$q3.Cells.Item(1, 5).Value = "股票总仓位"
$q3.Cells.Item(1, 6).Value = "仓位占比"
$q3.Cells.Item(1, 7).Value = "持有市值(亿元)"
$q3.Cells.Item(1, 8).Value = "仓位排名"

# Force columns B (fund code, e.g. "090019") and D:G (decimal-looking
# figures stored as text in the source data) to Text so leading zeros /
# exact string formatting survive instead of being coerced to numbers.
$q3.Range("B2:B5").NumberFormat = "@"
$q3.Range("D2:G5").NumberFormat = "@"

# Data rows
$q3.Cells.Item(2, 1).Value = 0
$q3.Cells.Item(2, 2).Value = "090019"
$q3.Cells.Item(2, 3).Value = "大成景恒混合A"
$q3.Cells.Item(2, 4).Value = "1.13"
$q3.Cells.Item(2, 5).Value = "93.98"
$q3.Cells.Item(2, 6).Value = "2.07"
$q3.Cells.Item(2, 7).Value = "0.0234"
$q3.Cells.Item(2, 8).Value = 1

$q3.Cells.Item(3, 1).Value = 1
$q3.Cells.Item(3, 2).Value = "009658"
$q3.Cells.Item(3, 3).Value = "汇丰晋信中小盘低波动策略股票A"
$q3.Cells.Item(3, 4).Value = "0.85"
$q3.Cells.Item(3, 5).Value = "90.14"
$q3.Cells.Item(3, 6).Value = "1.91"
$q3.Cells.Item(3, 7).Value = "0.0162"
$q3.Cells.Item(3, 8).Value = 10

$q3.Cells.Item(4, 1).Value = 2
$q3.Cells.Item(4, 2).Value = "006038"
$q3.Cells.Item(4, 3).Value = "大成景恒混合C"
$q3.Cells.Item(4, 4).Value = "0.45"
$q3.Cells.Item(4, 5).Value = "93.98"
$q3.Cells.Item(4, 6).Value = "2.07"
$q3.Cells.Item(4, 7).Value = "0.0093"
$q3.Cells.Item(4, 8).Value = 1

$q3.Cells.Item(5, 1).Value = 3
$q3.Cells.Item(5, 2).Value = "009775"
$q3.Cells.Item(5, 3).Value = "汇丰晋信中小盘低波动策略股票C"
$q3.Cells.Item(5, 4).Value = "0.04"
$q3.Cells.Item(5, 5).Value = "90.14"
$q3.Cells.Item(5, 6).Value = "1.91"
$q3.Cells.Item(5, 7).Value = "0.0008"
$q3.Cells.Item(5, 8).Value = 10

# Column A (running index) and H (rank) carry the same bordered/centered
# style the other quarter sheets use for column A.
$q3.Range("A2:A5").Style = $zongji.Range("A2").Style

# ---------------------------------------------------------------------
# 2) Insert a new row 2 into "总计" for the 2022-Q3 summary, and
#    renumber the running index in column A for the rows that shift down.
# ---------------------------------------------------------------------
$zongji.Rows.Item(2).Insert()

$zongji.Cells.Item(2, 1).Value = 0
$zongji.Cells.Item(2, 2).Value = "2022-Q3"
$zongji.Cells.Item(2, 3).Value = 4
$zongji.Cells.Item(2, 4).Value = 0.05

# Give the new A2 the same style as the other index cells in column A,
# and clear the stray formatting Insert() copied into B2:D2.
$zongji.Range("B2:D2").ClearFormats()
$zongji.Range("A2").Style = $zongji.Range("A3").Style

# The rows that shifted down (now rows 3-8) keep their original data but
# their running index (column A) needs to become 1,2,3,4,5,6.
$zongji.Cells.Item(3, 1).Value = 1
$zongji.Cells.Item(4, 1).Value = 2
$zongji.Cells.Item(5, 1).Value = 3
$zongji.Cells.Item(6, 1).Value = 4
$zongji.Cells.Item(7, 1).Value = 5
$zongji.Cells.Item(8, 1).Value = 6

Write-Output "2022-Q3 sheet + summary row added"
